$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns I and J, matching the style of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Add the data values for the new columns
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
